$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.315.96"
$ws.Range("E2").Value = "  -1.79%  "
$ws.Range("D3").Value = "2.596.24"
$ws.Range("E3").Value = "  -2.04%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.86%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.582"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.109"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.76"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.384"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.51%  "
$ws.Range("E12").Value = "  -0.75%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.47"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.43%  "
$ws.Range("D14").Value = "3.065.25"
$ws.Range("E14").Value = "  -2.00%  "
$ws.Range("D15").Value = "63.166.61"
$ws.Range("E15").Value = "  -1.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000154"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.26%  "
$ws.Range("D17").Value = "2.621.72"
$ws.Range("E17").Value = "  -0.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "343.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.63"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "565.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.161"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.63%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.02"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.58%  "
$ws.Range("D32").Value = "0.0₃0836"
$ws.Range("E32").Value = "  -4.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.27"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "166.11"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.410"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.39%  "
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.26"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.90"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.38%  "
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "166.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.96"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0578"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.21%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.08"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.97%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.628"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0246"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0958"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.80%  "
$ws.Range("D50").Value = "0.0₆0229"
$ws.Range("E50").Value = "  +14.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.179"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.96%  "
